# SANTA CRUZ DE LAS AGUAS.docx - update appointment schedule table
#
# The first data row (08:45 - 09:00 | (empty) | FLOR A FRUTO) is removed.
# The following two rows' time slots shift forward by one hour and fifteen
# minutes (09:00-09:15 -> 10:00-10:15, 09:15-09:30 -> 10:15-10:30), keeping
# their original buyer names (BOX BRAND, INTERLINK2AMERICAS).
# A new final row is appended with the slot 10:30 - 10:45 and buyer
# "FLOR A FRUTO" (the same buyer that used to occupy the first, now-removed,
# row).

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Remove the original first data row: "08:45 - 09:00" / "" / "FLOR A FRUTO"
$t.Rows.Item(2).Delete()

# Shift the time slots of the remaining two data rows forward
$t.Cell(2, 1).Range.Text = "10:00 - 10:15"
$t.Cell(3, 1).Range.Text = "10:15 - 10:30"

# Append a new row at the end of the table for the new slot
$t.Rows.Add() | Out-Null
$newRowIndex = $t.Rows.Count
$t.Cell($newRowIndex, 1).Range.Text = "10:30 - 10:45"
$t.Cell($newRowIndex, 3).Range.Text = "FLOR A FRUTO"
